# Updates the league base, swapping the match-detail data (columns B:AD)
# between pairs of rows, while keeping each row's sequential "id" (column A)
# fixed in place.
#
# Affected row pairs (1-based worksheet rows):
#   20 <-> 21
#   215 <-> 216
#   226 <-> 227
#   252 <-> 253

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($Row1, $Row2)

    $addr1 = "B" + $Row1 + ":AD" + $Row1
    $addr2 = "B" + $Row2 + ":AD" + $Row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $data1 = $range1.Value()
    $data2 = $range2.Value()

    $range1.Value = $data2
    $range2.Value = $data1
}

Swap-RowData 20 21
Swap-RowData 215 216
Swap-RowData 226 227
Swap-RowData 252 253
